$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Mars (planet)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "planet"
$ws.Range("C6").Value = "null"
$ws.Range("D6").Value = "Mars"
$ws.Range("E6").Value = "We will conquer the entire galaxy!!!"
$ws.Range("F6").Value = 999
$ws.Range("G6").Value = "null"
$ws.Range("H6").Value = 0.04
$ws.Range("I6").Value = "images/mars.jpg"

# Row 7: Norman (planet)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "planet"
$ws.Range("C7").Value = "null"
$ws.Range("D7").Value = "Norman"
$ws.Range("E7").Value = "MOOOOO!"
$ws.Range("F7").Value = 100000
$ws.Range("G7").Value = "null"
$ws.Range("H7").Value = 0.5
$ws.Range("I7").Value = "images/norman.jpg"

# Match style of existing data rows: columns A-F & I use font size 14 (no special number format);
# columns G & H use font size 14 with a 2-decimal number format.
$ws.Range("A6:F7").Font.Size = 14
$ws.Range("I6:I7").Font.Size = 14
$ws.Range("G6:H7").Font.Size = 14
$ws.Range("G6:H7").NumberFormat = "0.00"

# Update selection to match the diff (A6:XFD6 selected, active cell A6)
$ws.Range("A6:XFD6").Select()
